$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new tracked-time row (row 12): day / start / end / description.
$ws.Range("A12").Value = 43982
$ws.Range("B12").Value = 0.46527777777777773
$ws.Range("C12").Value = 0.82013888888888886
$ws.Range("E12").Value = "Google- Login- Lawyer design- register"

# Copy the "commit hash" cell format (small Consolas font, style index 3)
# from an existing populated cell (F2) onto F11, then set its value.
$ws.Range("F2").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value = "a59e8830d7d043c602947c1a73068a297c1412ef"

# Move the active selection, matching the saved view state.
$ws.Range("K14").Select() | Out-Null
